$wb = $excel.ActiveWorkbook

# --- Rename the first sheet ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Cadastro"

# --- Update login value (A2) from CaiqueOliveira to joaopedro ---
$ws1.Range("A2").Value = "joaopedro"

# --- Update selection on sheet1: was L1, now A2 ---
$ws1.Range("A2").Select()

# --- Add second sheet: "Pesquisa pagina inicial" (empty), placed after Cadastro ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Pesquisa pagina inicial"

# --- Add third sheet: "Pesquisa pela lupa" with content ---
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "Pesquisa pela lupa"
$ws3.Range("A1").Value = "HP PAVILION 15Z TOUCH LAPTOP"
$ws3.Range("A1").Font.Underline = 2

# --- Make the third sheet the active/selected tab ---
$ws3.Activate()
$ws3.Select()
